$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '58.166.03'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -0.42%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '2.593.39'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -0.92%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.08%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '521.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +0.11%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '143.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +0.57%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.23%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.567'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -0.18%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '2.614.35'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.52%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '6.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -1.50%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -1.44%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '0.341'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +1.96%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -0.04%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '3.051.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -0.88%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '58.124.62'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -0.45%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '20.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -2.94%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = "'" + '2.629.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -0.82%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = "'" + '0.0000134'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -1.43%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '339.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +0.78%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '4.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -1.73%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '10.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -1.54%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '6.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +2.68%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -0.08%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '65.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +0.57%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '0.168'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +0.96%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '0.404'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -2.26%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = "'" + '2.723.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -0.29%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = "'" + '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +0.01%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -1.14%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  -5.91%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -0.10%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '6.11'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -6.24%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -0.49%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '18.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +0.06%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '149.69'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -0.23%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -2.20%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -4.38%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '0.857'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -4.10%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '0.859'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +1.27%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = "'" + '36.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -0.53%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = "'" + '1.46'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +1.99%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '3.53'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -2.47%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '0.996'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -0.29%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '0.607'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +0.86%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '270.48'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +1.07%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '10.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +0.19%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '0.0954'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -1.63%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '18.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -1.93%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '0.0521'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -1.61%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'" + '4.68'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +2.43%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = "'" + '1.965.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -3.20%  '
$ws.Range('E51').Style = 'Normal'
